# Add two new data rows (Tp 5 and 6) to the CRM accuracy tracking sheet.
# This mirrors extending the table on Sheet1 with one more reading
# (row 42), reusing the "Opened CRM (8/30/2019)" batch label already
# used by the previous row (41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell's formatting (style) from the row above so the new
# date cell keeps the same date number format instead of picking up a
# brand-new style entry.
$ws.Range("A41").Copy($ws.Range("A42"))

# New row of data.
$ws.Range("A42").Value = 43712
$ws.Range("B42").Value = 2218.9034025771698
$ws.Range("C42").Value = 2207.0300000000002
$ws.Range("D42").Formula = "=100*(B42-C42)/C42"
$ws.Range("E42").Value = 169
$ws.Range("F42").Value = "Opened CRM (8/30/2019)"

# Match the selection left behind in the saved file.
[void]$ws.Range("F42").Select()
